# "diálogos madre y chatarrero"
# Replace the placeholder Miranda ("Inicio") and Chatarero dialogue rows with
# the real Spanish dialogue lines, and renumber / relabel accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Miranda's first conversation ("Miranda1") ---
$ws.Range("A9").Value  = "Miranda1"
$ws.Range("B9").Value  = 0
$ws.Range("C9").Value  = 1
$ws.Range("D9").Value  = "Miranda"
$ws.Range("E9").Value  = "Miranda"
$ws.Range("F9").Value  = "Hola hijo, buenos días ¿Dónde estabas?"

$ws.Range("A10").Value = "Miranda1"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = "Miranda"
$ws.Range("E10").Value = "Miranda"
$ws.Range("F10").Value = "Bueno es igual, necesito que me hagas un favor."

$ws.Range("A11").Value = "Miranda1"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = "Miranda"
$ws.Range("E11").Value = "Miranda"
$ws.Range("F11").Value = "Por favor, ve a la ferreteria del prueblo a por un poco de  aceite."

$ws.Range("A12").Value = "Miranda1"
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = "Miranda"
$ws.Range("E12").Value = "Miranda"
$ws.Range("F12").Value = "Esta nueva enfermedad me está afectando más de lo que esperaba."

$ws.Range("A13").Value = "Miranda1"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 5
$ws.Range("D13").Value = "Miranda"
$ws.Range("E13").Value = "Miranda"
$ws.Range("F13").Value = "Muchas gracias hijo."

# --- Miranda's second conversation ("Miranda2") ---
$ws.Range("A14").Value = "Miranda2"
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 6
$ws.Range("D14").Value = "Miranda"
$ws.Range("E14").Value = "Miranda"
$ws.Range("F14").Value = "Barry, debes irte de Ferder, el lugar ya no es seguro."

$ws.Range("A15").Value = "Miranda2"
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = "Miranda"
$ws.Range("E15").Value = "Miranda"
$ws.Range("F15").Value = "La enfermedad nos esta pasando factura a todos, y no hay nadie que nos pueda ayudar."

$ws.Range("A16").Value = "Miranda2"
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = "Miranda"
$ws.Range("E16").Value = "Miranda"
$ws.Range("F16").Value = "Si te vas lejos puede que encuentres a alguien que sepa como ayudarnos."

$ws.Range("A17").Value = "Miranda2"
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = "Miranda"
$ws.Range("E17").Value = "Miranda"
$ws.Range("F17").Value = "Toma el camino del bosque, quizá el granjero Hamilton te pueda ayudar."

# --- Chatarero's conversation ---
$ws.Range("E18").Value = "Chatarero"
$ws.Range("F18").Value = "Bienvenido pequeño, ¿Qué necesitas?"

$ws.Range("E19").Value = "Chatarero"
$ws.Range("F19").Value = "O vaya, así que un poco de aceite para tu mamá."

$ws.Range("F20").Value = "Espera un segundo, enseguida te lo doy."

$ws.Range("F21").Value = "… "

$ws.Range("F22").Value = "… "

$ws.Range("E23").Value = "Chatarero"
$ws.Range("F23").Value = "… "

$ws.Range("E24").Value = "Chatarero"
$ws.Range("F24").Value = "… "

$ws.Range("F25").Value = "Perfecto, aquí lo tienes."

$ws.Range("F26").Value = "Espero que esta nueva enfermedad se acabe ya, o esto no afectará a todos."

# Update the selection to match the author's last-saved cursor position.
[void]$ws.Range("I8").Select()
